$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column D ("Allergen"), which will
# hold the ingredient's E number. Existing columns D.. shift right to E.. .
$ws.Columns("D:D").Insert()

# New header for the inserted column.
$ws.Range("D1").Value = "Enumber"
$ws.Range("D1").Font.Bold = $true

# Match the column width style used for other text columns instead of the
# bestFit numeric width that Excel would otherwise guess.
$ws.Columns("D:D").ColumnWidth = 20.08

# Populate E numbers for the relevant ingredient rows.
$ws.Range("D3").Value = 300
$ws.Range("D5").Value = 150
$ws.Range("D6").Value = 290
$ws.Range("D7").Value = 469
$ws.Range("D8").Value = 330
$ws.Range("D16").Value = 297
$ws.Range("D18").Value = 414
$ws.Range("D19").Value = 270
$ws.Range("D20").Value = 296
$ws.Range("D21").Value = 353
$ws.Range("D27").Value = 202
$ws.Range("D32").Value = 334

# Reset the saved cursor position back to the default top-left cell.
$ws.Range("A1").Select()
